# SACR Complete Changelog - Edition 2: add "Changelog #8: SACR R4.2" section
# (Updated Rig GitHub to R4.2)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text of a Range (excluding its trailing paragraph
# mark) with $text. Re-setting identical text is a silent no-op in this
# runtime, so stage a throw-away placeholder first to force the rewrite
# (this is what collapses the multiple same-text runs into one run, as
# happens in the diff for the "Changelog #8: SACR R4.1.2 LTS" heading).
# ---------------------------------------------------------------------------
function Set-ParaText($para, [string]$text) {
    $rng = $para.Range
    $body = $d.Range($rng.Start, $rng.End - 1)
    $body.Text = "@@__TMP__@@"
    $para2 = $d.Paragraphs.Item($para.Index)
    $rng2 = $para2.Range
    $body2 = $d.Range($rng2.Start, $rng2.End - 1)
    $body2.Text = $text
}

# ---------------------------------------------------------------------------
# 1) Collapse the 4 runs of the existing "Changelog #8: SACR R4.1.2 LTS"
#    Heading 1 into a single run (no visible text change).
# ---------------------------------------------------------------------------
$oldHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    # Paragraph.Range.Text includes the trailing paragraph mark (\r), so
    # compare with StartsWith rather than an exact match.
    if ($cand.Range.Text.StartsWith("Changelog #8: SACR R4.1.2 LTS")) {
        $oldHeading = $cand
        break
    }
}
Set-ParaText $oldHeading "Changelog #8: SACR R4.1.2 LTS"

# ---------------------------------------------------------------------------
# 2) Replace the final paragraph (the one holding the inline picture) with
#    a brand-new "Changelog #8: SACR R4.2" section.
# ---------------------------------------------------------------------------

# Find the paragraph that carries the <w:drawing> (last paragraph, no text).
$picPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$picPara.Range.Delete()

# Templates to clone paragraph-mark / run formatting (font lang tag) from,
# matched by style + list level, so the new paragraphs serialize with the
# same implicit formatting as the rest of the document.
$tplHeading1 = $null
$tplHeading2 = $null
$tplList1 = $null
$tplList2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $styleName = $cand.Style.NameLocal
    if ($styleName -eq "Heading 1" -and $tplHeading1 -eq $null) {
        $tplHeading1 = $cand
    }
    if ($styleName -eq "Heading 2" -and $tplHeading2 -eq $null) {
        $tplHeading2 = $cand
    }
    if ($styleName -eq "List Paragraph") {
        $lvl = $cand.Range.ListFormat.ListLevelNumber
        if ($lvl -eq 2 -and $tplList1 -eq $null) {
            $tplList1 = $cand
        }
        if ($lvl -eq 3 -and $tplList2 -eq $null) {
            $tplList2 = $cand
        }
    }
}

# Appends a new paragraph at the end of the document with the given Word
# style name, outline/list level (0 = Heading1 "section" level with no
# bullet, 1 = top bullet level (ilvl 0), 2 = nested bullet (ilvl 1),
# 3 = doubly-nested bullet (ilvl 2)) and text.
function Add-Para([string]$styleName, [int]$level, [string]$text) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newRange = $newPara.Range
    $startPos = $newRange.Start

    switch ($level) {
        0 { $tpl = $tplHeading1 }
        1 { $tpl = $tplHeading2 }
        2 { $tpl = $tplList1 }
        3 { $tpl = $tplList2 }
    }

    $newRange.FormattedText = $tpl.Range.FormattedText
    $newPara.Style = $styleName

    $tplLen = $tpl.Range.End - $tpl.Range.Start
    $body = $d.Range($startPos, $startPos + $tplLen)
    $body.Text = $text

    $finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    if ($level -ge 1) {
        $finalPara.Range.ListFormat.ListLevelNumber = $level
    }
}

Add-Para "Heading 1" 0 "Changelog #8: SACR R4.2"

Add-Para "Heading 2" 1 "Fixes"
Add-Para "List Paragraph" 2 "Fixed Parenting of Extrude Meshes"
Add-Para "List Paragraph" 2 "Disabled Denoiser by Default"
Add-Para "List Paragraph" 3 "For older systems that don’t support SSE4.2"
Add-Para "List Paragraph" 2 "Re-enabled Extra’s and Origin’s in viewport display settings"
Add-Para "List Paragraph" 2 "Leg Lattice Alignment"
Add-Para "List Paragraph" 2 "Fixed Eye Gradient UV’s"
Add-Para "List Paragraph" 2 "Teeth Alignment"
Add-Para "List Paragraph" 2 "Fixed Emission Controls for Eyewhites not working"
Add-Para "List Paragraph" 2 "Fixed Subdivision Surface"
Add-Para "List Paragraph" 3 "Removed Old Depricated Drivers and optimized"
Add-Para "List Paragraph" 2 "Optimized Memory Usage"

Add-Para "Heading 2" 1 "Changes"
Add-Para "List Paragraph" 2 "Edited Default Material Values"
Add-Para "List Paragraph" 2 "Changed Default Textures"
Add-Para "List Paragraph" 3 "Skin changed to Medieval (No Face)"
Add-Para "List Paragraph" 3 "HD Eye changed to Medieval (with Face) "
Add-Para "List Paragraph" 3 "HD Eye White"
Add-Para "List Paragraph" 2 "Boneshape Re-design"
Add-Para "List Paragraph" 2 "Dropped LTS Tag (Final R4 build)"

Add-Para "Heading 2" 1 "Additions"
Add-Para "List Paragraph" 2 "Easy Parent Objects"
Add-Para "List Paragraph" 3 "Compatible with head and wrists"

Add-Para "Heading 2" 1 "Removed"
Add-Para "List Paragraph" 2 "Leg IK/FK Hybrid System"
Add-Para "List Paragraph" 3 "Reverted to traditional switch"

Write-Output ("Done. Paragraph count=" + $d.Paragraphs.Count)
